$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")

# --- Update header row (D2: SortOrder header stays same text, but the
#     shared-string it pointed at moves because other strings are removed) ---
# (No value change needed for header text itself - already "SortOrder";
#  the shared string index shift happens automatically on save.)

# --- Row 3 (Id 1000, the root/site node): rename Name from "example.com" to "Site" ---
$ws.Range("F3").Value = "Site"

# --- Row 4 (Id 1001): rename Name/listHeading from "About us" to "Item 1"
#     and clear the listSummary ("Read all about us" html) ---
$ws.Range("F4").Value = "Item 1"
$ws.Range("H4").Value = "Item 1"
$ws.Range("I4").ClearContents()

# --- Row 5 (Id 1002): rename Name/listHeading from "Contact us" to "Item 2"
#     and clear the listSummary ("Contact us here" html) ---
$ws.Range("F5").Value = "Item 2"
$ws.Range("H5").Value = "Item 2"
$ws.Range("I5").ClearContents()

# --- Add new stub rows 6-9 (children of 1001 and 1002) ---
$ws.Range("B6").Value = 1003
$ws.Range("C6").Value = 1001
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 1057
$ws.Range("F6").Value = "Item 1.1"
$ws.Range("G6").Value = 1056
$ws.Range("H6").Value = "Item 1.1"

$ws.Range("B7").Value = 1004
$ws.Range("C7").Value = 1001
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 1057
$ws.Range("F7").Value = "Item 1.2"
$ws.Range("G7").Value = 1056
$ws.Range("H7").Value = "Item 1.2"

$ws.Range("B8").Value = 1005
$ws.Range("C8").Value = 1002
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 1057
$ws.Range("F8").Value = "Item 2.1"
$ws.Range("G8").Value = 1056
$ws.Range("H8").Value = "Item 2.1"

$ws.Range("B9").Value = 1006
$ws.Range("C9").Value = 1002
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 1057
$ws.Range("F9").Value = "Item 2.2"
$ws.Range("G9").Value = 1056
$ws.Range("H9").Value = "Item 2.2"

# --- Resize Table1 to cover the newly added rows ---
$lo = $ws.ListObjects.Item("Table1")
$lo.Resize($ws.Range("B2:I9"))

# --- Fix up sheet view/selection: Input becomes the active/selected sheet,
#     with the JSON sheet no longer active, and selection moved to I9 ---
$ws.Activate()
$ws.Range("I9").Select()
